# 書籍一覧.xlsx – convert the 出版日 (publish date) column from real dates
# to plain text strings formatted as "yyyy/mm/dd", and give 価格 (price)
# an explicit #,##0 number format. This mirrors what happens when a user
# selects the 出版日 column of the table and switches it to Text format,
# then retypes the dates as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column keeps its #,##0 number formatting.
$ws.Range("D2:D4").NumberFormat = "#,##0"

# Publish-date column becomes Text (@) formatted ...
$ws.Range("C2:C4").NumberFormat = "@"

# ... and the dates are re-entered as literal text (not real dates).
$ws.Range("C2").Value = "2022/03/10"
$ws.Range("C3").Value = "2021/11/12"
$ws.Range("C4").Value = "2021/07/15"

# Column C is widened to fit the new text values.
$ws.Columns("C").AutoFit() | Out-Null

# Move the active selection off the table, onto G3.
$ws.Range("G3").Select() | Out-Null
